$wb = $excel.ActiveWorkbook
$planSheet = $wb.Worksheets.Item("Plan Details")

# New rows for the "Legs" section (Squats / Leg Press) describing the
# workout creation API parse output: Sets / Rep-Range / Exercise Description
$planSheet.Range("H11").Value = 3
$planSheet.Range("I11").Value = "5 to 9"
$planSheet.Range("J11").Value = "parallel to ground"

$planSheet.Range("H12").Value = 3
$planSheet.Range("I12").Value = "5 to 9"
$planSheet.Range("J12").Value = "pause at lowest"

# Bench Press description tweak (trailing space variant added as a new entry)
$planSheet.Range("J2").Value = "slow negatives "

# Make "Plan Details" the active sheet / tab, with J2 selected
$planSheet.Activate()
$planSheet.Range("J2").Select() | Out-Null
